$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the top of the time series
# for this market/product. Insert a new row at 88 (shifting the existing
# rows 88-137 down to 89-138) and populate it with the new week's data.
$ws.Rows.Item(88).Insert()

$ws.Range("A88").Value = 11
$ws.Range("B88").Value = 'Vega Monumental Concepción'
$ws.Range("C88").Value = 'Bíobío'
$ws.Range("D88").Value = 44596
$ws.Range("E88").Value = 8
$ws.Range("F88").Value = 100112003
$ws.Range("G88").Value = 'Ajo'
$ws.Range("H88").Value = 'Chino'
$ws.Range("I88").Value = 'Primera'
$ws.Range("J88").Value = 400
$ws.Range("K88").Value = 17000
$ws.Range("L88").Value = 18000
$ws.Range("M88").Value = 17500
$ws.Range("N88").Value = '$/caja 10 kilos'
$ws.Range("O88").Value = 'China'
$ws.Range("P88").Value = 1750
$ws.Range("Q88").Value = 10
$ws.Range("R88").Value = 'Hortaliza'
